$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# --- Row 2 (CasesTab) ---
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Value = "CasesTab"

$q1 = 'MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_type = ''Aligned RNA reads file'' 
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity'
$ws.Range("B2").Value = $q1
$ws.Range("B2").WrapText = $true

$q2 = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_type = ''Aligned RNA reads file''
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'
$ws.Range("C2").Value = $q2
$ws.Range("C2").WrapText = $true

$ws.Range("D2").Value = "TC02_Trials_Filter_AssocFileType-AlignedRNA_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC02_Trials_Filter_AssocFileType-AlignedRNA_WebData.xlsx"

$ws.Rows(2).RowHeight = 188.5

# --- Row 3 (FilesTab) ---
$ws.Range("A3").Value = "FilesTab"

$q3 = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE f.file_type = ''Aligned RNA reads file''
 WITH
    f, parent, c, a, ct,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`'
$ws.Range("B3").Value = $q3
$ws.Range("B3").WrapText = $true

$ws.Range("C3").Value = $q2
$ws.Range("C3").WrapText = $true

$ws.Range("D3").Value = "TC02_Trials_Filter_AssocFileType-AlignedRNA_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC02_Trials_Filter_AssocFileType-AlignedRNA_WebData.xlsx"

$ws.Rows(3).RowHeight = 409.5

# --- Column widths ---
$ws.Columns(1).ColumnWidth = 8
$ws.Columns(2).ColumnWidth = 75
$ws.Columns(3).ColumnWidth = 75
$ws.Columns(4).ColumnWidth = 69.5
$ws.Columns(5).ColumnWidth = 27.66666666666667

# --- Selection / view ---
$ws.Range("C3").Select()
